$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (E1) onto the new header cell (F1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header text
$ws.Range("F1").Value = "Modelo"

# Updated metric values
$ws.Range("B2").Value = 0.2486823901655659
$ws.Range("C2").Value = 0.9951394964303211
$ws.Range("D2").Value = 0.3922038828384552

# New model name column
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=3, n_estimators=50))])"
